$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Type the same values a user would have entered by hand into column W.
$ws.Range("W1").Value = 1
$ws.Range("W3").Value = 1
$ws.Range("W5").Value = 1

# After typing into W5 and pressing Enter, the active cell moves to W6 --
# leave the selection there, matching the saved sheet view.
[void]$ws.Range("W6").Select()
